$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 378 ("Región de Arica y
# Parinacota", 2023-08-04), pushing the existing rows 378-423 down to
# 379-424 (last row becomes 424).
$ws.Rows.Item(378).Insert()

$ws.Cells.Item(378, 1).Value = 10
$ws.Cells.Item(378, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(378, 3).Value = "La Araucanía"
$ws.Cells.Item(378, 4).Value = 45142
$ws.Cells.Item(378, 5).Value = 9
$ws.Cells.Item(378, 6).Value = 100112052
$ws.Cells.Item(378, 7).Value = "Albahaca"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 80
$ws.Cells.Item(378, 11).Value = 6000
$ws.Cells.Item(378, 12).Value = 6000
$ws.Cells.Item(378, 13).Value = 6000
$ws.Cells.Item(378, 14).Value = "$/paquete"
$ws.Cells.Item(378, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(378, 16).Value = 6000
$ws.Cells.Item(378, 17).Value = 1
$ws.Cells.Item(378, 18).Value = "Hortaliza"
